$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OutOfStock")

$row = 13
$ws.Cells.Item($row, 1).Value = "x1"
$ws.Cells.Item($row, 2).Value = "Screwdriver"
$ws.Cells.Item($row, 3).Value = "2025-06-14T03:23:22.852Z"
$ws.Cells.Item($row, 4).Value = "14/6/2025, 8:53:22 am"
